# Updates the cryptos price/volume snapshot (Price column D, Volume(1h) column E)
# and fixes the row order / data for two swapped coin rows (34/35 and 43/44).
# Price values that would otherwise be auto-coerced to a number by Excel
# (e.g. "214.71") are written with a leading "'" (quote-prefix) so they stay
# text, matching the source data which stores every Price cell as a string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.025.45'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '1.633.92'
$ws.Range('E3').Value = '  -0.74%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '''214.71'
$ws.Range('E5').Value = '  -0.88%  '
$ws.Range('E6').Value = '  -1.15%  '
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('D8').Value = '''0.251'
$ws.Range('E8').Value = '  -2.28%  '
$ws.Range('E9').Value = '  -3.09%  '
$ws.Range('D10').Value = '''18.40'
$ws.Range('D11').Value = '''0.0791'
$ws.Range('E11').Value = '  -0.61%  '
$ws.Range('D12').Value = '1.860.35'
$ws.Range('E12').Value = '  -0.79%  '
$ws.Range('D13').Value = '1.633.57'
$ws.Range('E13').Value = '  -4.89%  '
$ws.Range('E14').Value = '  -2.69%  '
$ws.Range('D15').Value = '''0.526'
$ws.Range('E15').Value = '  -3.36%  '
$ws.Range('D16').Value = '26.000.06'
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('D17').Value = '0.0₃0742'
$ws.Range('E17').Value = '  -2.96%  '
$ws.Range('D18').Value = '''61.43'
$ws.Range('E18').Value = '  -3.16%  '
$ws.Range('E19').Value = '  +0.20%  '
$ws.Range('D20').Value = '''190.78'
$ws.Range('E20').Value = '  -2.49%  '
$ws.Range('D21').Value = '''4.25'
$ws.Range('E21').Value = '  -2.43%  '
$ws.Range('E22').Value = '  -3.05%  '
$ws.Range('D23').Value = '''6.09'
$ws.Range('E23').Value = '  -2.26%  '
$ws.Range('D24').Value = '''0.132'
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range('D25').Value = '''1.79'
$ws.Range('E25').Value = '  -1.01%  '
$ws.Range('D26').Value = '''143.64'
$ws.Range('E26').Value = '  -0.17%  '
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('E28').Value = '  -1.87%  '
$ws.Range('D29').Value = '''15.17'
$ws.Range('E29').Value = '  -2.46%  '
$ws.Range('E30').Value = '  -1.40%  '
$ws.Range('E31').Value = '  -3.24%  '
$ws.Range('E32').Value = '  -4.40%  '
$ws.Range('D33').Value = '''3.12'
$ws.Range('E33').Value = '  -5.42%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '''2.40'
$ws.Range('E34').Value = '  -2.28%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = '''1.49'
$ws.Range('E35').Value = '  -3.52%  '
$ws.Range('D36').Value = '1.131.43'
$ws.Range('E36').Value = '  -0.21%  '
$ws.Range('D37').Value = '''0.862'
$ws.Range('E37').Value = '  -5.03%  '
$ws.Range('E38').Value = '  -1.08%  '
$ws.Range('E39').Value = '  -4.48%  '
$ws.Range('E40').Value = '  -1.62%  '
$ws.Range('E41').Value = '  -1.08%  '
$ws.Range('E42').Value = '  -2.97%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '''5.23'
$ws.Range('E43').Value = '  -5.04%  '
$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = '1.769.40'
$ws.Range('E44').Value = '  -0.86%  '
$ws.Range('E45').Value = '  -2.05%  '
$ws.Range('D46').Value = '''54.83'
$ws.Range('E46').Value = '  -3.39%  '
$ws.Range('D47').Value = '''0.0527'
$ws.Range('E47').Value = '  +0.05%  '
$ws.Range('E48').Value = '  +1.23%  '
$ws.Range('E49').Value = '  -0.22%  '
$ws.Range('E50').Value = '  +0.35%  '
$ws.Range('D51').Value = '''7.50'
$ws.Range('E51').Value = '  -3.66%  '
